# "9th Stab - Cosmetic Changes"
# Column C duplicated column B's data (with a couple of distinct
# "latest rating" values and one highlighted cell). Remove the redundant
# column: fold any distinct values from C into B, preserve the highlight
# fill on the one cell that carried it, then delete column C so the
# sheet collapses back down to two columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 27 }

for ($r = 1; $r -le $lastRow; $r++) {
    $src = $ws.Cells.Item($r, 3)
    $val = $src.Value2
    if ($val -ne $null) {
        $dst = $ws.Cells.Item($r, 2)
        $dst.Value2 = $val

        $colorIndex = $src.Interior.ColorIndex
        if ($colorIndex -ne -4142) {
            $dst.Interior.ColorIndex = $colorIndex
        }
    }
}

$ws.Columns.Item(3).Delete() | Out-Null

$ws.Range("D8").Select() | Out-Null
